$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

# Column A holds a date-looking string ("2025-10-22"); force Text format
# first so Excel doesn't silently convert it into a date serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-10-22"
$ws.Cells.Item($row, 2).Value = "15:21:44"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,801.9975"
